$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arduino")

# Row 26: keep existing "Наклейки" label, add price
$ws.Range("G26").Value = 333

# Row 27: was "Часы реального времени" -> becomes "SD карта", add price
$ws.Range("B27").Value = "SD карта"
$ws.Range("G27").Value = 500

# Row 28: was "Wf модуль" -> becomes "Вилка 32А " with article + price
$ws.Range("B28").Value = "Вилка 32А "
$ws.Range("G28").Value = 375

# Row 29 (new): "Розетка силовая" with article + price
$ws.Range("B29").Value = "Розетка силовая"
$ws.Range("C29").Value = "223 2P+N32A"
$ws.Range("G29").Value = 251

# Now set the article for row 28 (ensures shared string insertion order
# matches: SD карта, Вилка 32А , Розетка силовая, 223 2P+N32A, 523 2P + N 32А)
$ws.Range("C28").Value = "523 2P + N 32А"

# Rows 31/32 (new): push the displaced "Часы реального времени" / "Wf модуль"
# entries further down the list, leaving row 30 blank
$ws.Range("B31").Value = "Часы реального времени"
$ws.Range("B32").Value = "Wf модуль"

# Extend the total formula to cover the newly added rows
$ws.Range("G1").Formula = "=SUM(G4:G39)"

# Restore the selection shown in the saved workbook
$ws.Range("J29").Select()
